# Apply the "Update Excel template Samit - new subject therapeutic tutoring, Sir rakaz" edit.
$wb = $excel.ActiveWorkbook

$wsTables = $wb.Worksheets.Item("טבלאות")

# Insert the new "חונכות טיפולית" (therapeutic tutoring) subject at the top of the
# subject list (column A), pushing the existing entries down by one row, and also
# append it again as a new list entry at the end of the currently-used range.
$wsTables.Range("A3").Value = "חונכות טיפולית"
$wsTables.Range("A4").Value = "שיעורי עזר"
$wsTables.Range("A5").Value = "הוראה מתקנת"
$wsTables.Range("A6").Value = "תרגום"
$wsTables.Range("A7").Value = "חונכות טיפולית"

# Add the new "שיר" (Sir) coordinator name to the coordinators list (column C).
$wsTables.Range("C14").Value = "שיר"

# Extend the "רכז תלמיד" (student coordinator) drop-down on the main sheet so it
# covers the newly added list entry.
$wsMain = $wb.Worksheets.Item("ראשי")
$wsMain.Range("L2").Validation.Delete()
$wsMain.Range("L2").Validation.Add(3, 1, 1, "=טבלאות!`$C`$2:`$C`$14")

# Leave the cursor selection on L2 / C3:C14 the way the authored workbook does.
$wsMain.Range("L2").Select()
$wsTables.Range("C3:C14").Select()
$wsTables.Range("C3").Activate()
